# Auto-generated edit script applying numeric corrections to Sheets/Omega_Profits.xlsx
# as described by the commit diff (scheduled runner market data refresh).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 198.125
$ws.Range("I6").Value = 205
$ws.Range("K6").Value = 615
$ws.Range("M6").Value = -503
# Row 11
$ws.Range("H11").Value = 499.8
$ws.Range("I11").Value = 499.8
$ws.Range("K11").Value = 499.8
$ws.Range("M11").Value = -359.8
# Row 19
$ws.Range("H19").Value = 2809.7144
$ws.Range("I19").Value = 2529.3333
$ws.Range("K19").Value = 2529.3333
$ws.Range("M19").Value = -2354.3333
# Row 40
$ws.Range("H40").Value = 7215.0454
$ws.Range("J40").Value = 10749.846
$ws.Range("L40").Value = 10749.846
$ws.Range("N40").Value = -11099.846
# Row 43
$ws.Range("H43").Value = 4833.154
$ws.Range("I43").Value = 4673.25
$ws.Range("J43").Value = 4904.222
$ws.Range("K43").Value = 4673.25
$ws.Range("L43").Value = 4904.222
$ws.Range("M43").Value = -4604.25
$ws.Range("N43").Value = -5042.222
# Row 107
$ws.Range("H107").Value = 2227.3076
$ws.Range("I107").Value = 1995.8
$ws.Range("K107").Value = 1995.8
$ws.Range("M107").Value = -75.79999999999995
# Row 132
$ws.Range("H132").Value = 2591.257
$ws.Range("I132").Value = 2571.9033
$ws.Range("K132").Value = 7715.7099
$ws.Range("M132").Value = -5185.7099
# Row 137
$ws.Range("H137").Value = 4998
$ws.Range("I137").Value = 4997
$ws.Range("K137").Value = 14991
$ws.Range("M137").Value = -12441

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 10192.667
$ws.Range("I61").Value = 7093.643
$ws.Range("K61").Value = 7093.643
$ws.Range("M61").Value = -6881.643
# Row 63
$ws.Range("H63").Value = 5011.8096
$ws.Range("I63").Value = 5011.8096
$ws.Range("K63").Value = 5011.8096
$ws.Range("M63").Value = -4325.8096
# Row 66
$ws.Range("H66").Value = 5011.8096
$ws.Range("I66").Value = 5011.8096
$ws.Range("K66").Value = 25059.048
$ws.Range("M66").Value = -21627.048
# Row 74
$ws.Range("H74").Value = 1950.1212
$ws.Range("J74").Value = 2268.6875
$ws.Range("L74").Value = 2268.6875
$ws.Range("N74").Value = -4016.6875
# Row 77
$ws.Range("H77").Value = 1950.1212
$ws.Range("J77").Value = 2268.6875
$ws.Range("L77").Value = 11343.4375
$ws.Range("N77").Value = -20079.4375
# Row 102
$ws.Range("H102").Value = 1548.5264
$ws.Range("I102").Value = 1534.7142
$ws.Range("K102").Value = 1534.7142
$ws.Range("M102").Value = 87.28580000000011
# Row 133
$ws.Range("H133").Value = 62847.4
$ws.Range("J133").Value = 58500
$ws.Range("L133").Value = 58500
$ws.Range("N133").Value = -63560
# Row 136
$ws.Range("H136").Value = 10192.667
$ws.Range("I136").Value = 7093.643
$ws.Range("K136").Value = 21280.929
$ws.Range("M136").Value = -18730.929

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 13890396
$ws.Range("I94").Value = 13890396
$ws.Range("K94").Value = 13890396
$ws.Range("M94").Value = -13889945
# Row 134
$ws.Range("H134").Value = 2389.75
$ws.Range("I134").Value = 2362
$ws.Range("J134").Value = 2436
$ws.Range("K134").Value = 7086
$ws.Range("L134").Value = 7308
$ws.Range("M134").Value = -4551
$ws.Range("N134").Value = -12378

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 14
$ws.Range("H14").Value = 20000
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 20000
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 20000
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -20340
# Row 31
$ws.Range("H31").Value = 8972.909
$ws.Range("I31").Value = 16579.8
$ws.Range("K31").Value = 16579.8
$ws.Range("M31").Value = -16284.8
# Row 34
$ws.Range("H34").Value = 8972.909
$ws.Range("I34").Value = 16579.8
$ws.Range("K34").Value = 16579.8
$ws.Range("M34").Value = -16377.8
# Row 41
$ws.Range("H41").Value = 43000
$ws.Range("J41").Value = 43000
$ws.Range("L41").Value = 43000
$ws.Range("N41").Value = -43856
# Row 58
$ws.Range("H58").Value = 5020
$ws.Range("I58").Value = 4300.3335
$ws.Range("K58").Value = 4300.3335
$ws.Range("M58").Value = -4097.3335
# Row 70
$ws.Range("H70").Value = 56999.5
$ws.Range("J70").Value = 56999.5
$ws.Range("L70").Value = 56999.5
$ws.Range("N70").Value = -57629.5
# Row 73
$ws.Range("H73").Value = 56999.5
$ws.Range("J73").Value = 56999.5
$ws.Range("L73").Value = 56999.5
$ws.Range("N73").Value = -59183.5
# Row 80
$ws.Range("H80").Value = 44127.5
$ws.Range("J80").Value = 44127.5
$ws.Range("L80").Value = 44127.5
$ws.Range("N80").Value = -46373.5
# Row 81
$ws.Range("H81").Value = 38327.5
$ws.Range("J81").Value = 38327.5
$ws.Range("L81").Value = 38327.5
$ws.Range("N81").Value = -40323.5
# Row 83
$ws.Range("H83").Value = 44127.5
$ws.Range("J83").Value = 44127.5
$ws.Range("L83").Value = 132382.5
$ws.Range("N83").Value = -143614.5
# Row 84
$ws.Range("H84").Value = 38327.5
$ws.Range("J84").Value = 38327.5
$ws.Range("L84").Value = 114982.5
$ws.Range("N84").Value = -124966.5
# Row 134
$ws.Range("H134").Value = 4491.84
$ws.Range("I134").Value = 4014.8572
$ws.Range("K134").Value = 12044.5716
$ws.Range("M134").Value = -9509.571599999999
# Row 136
$ws.Range("H136").Value = 5020
$ws.Range("I136").Value = 4300.3335
$ws.Range("K136").Value = 12901.0005
$ws.Range("M136").Value = -10351.0005

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 32
$ws.Range("H32").Value = 2183.1667
$ws.Range("I32").Value = 4250
$ws.Range("K32").Value = 12750
$ws.Range("M32").Value = -12467
# Row 128
$ws.Range("H128").Value = 165513
$ws.Range("I128").Value = 165513
$ws.Range("K128").Value = 496539
$ws.Range("M128").Value = -491559

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 4733537
$ws.Range("I3").Value = 3638947
$ws.Range("J3").Value = 6375422
$ws.Range("K3").Value = 3638947
$ws.Range("L3").Value = 6375422
$ws.Range("M3").Value = -3638831
$ws.Range("N3").Value = -6375654
# Row 4
$ws.Range("H4").Value = 500
$ws.Range("J4").Value = 500
$ws.Range("L4").Value = 500
$ws.Range("N4").Value = -724
# Row 15
$ws.Range("H15").Value = 20000
$ws.Range("J15").Value = 20000
$ws.Range("L15").Value = 20000
$ws.Range("N15").Value = -20576
# Row 43
$ws.Range("H43").Value = 1379.6
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
# Row 81
$ws.Range("H81").Value = 20000
$ws.Range("J81").Value = 20000
$ws.Range("L81").Value = 20000
$ws.Range("N81").Value = -21996
# Row 84
$ws.Range("H84").Value = 20000
$ws.Range("J84").Value = 20000
$ws.Range("L84").Value = 60000
$ws.Range("N84").Value = -69984
# Row 132
$ws.Range("H132").Value = 3037.1738
$ws.Range("I132").Value = 3382.6
$ws.Range("J132").Value = 2389.5
$ws.Range("K132").Value = 10147.8
$ws.Range("L132").Value = 7168.5
$ws.Range("M132").Value = -7617.799999999999
$ws.Range("N132").Value = -12228.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 24256.75
$ws.Range("I7").Value = 19426.584
$ws.Range("K7").Value = 19426.584
$ws.Range("M7").Value = -19314.584
# Row 41
$ws.Range("H41").Value = 15000
$ws.Range("I41").Value = 15000
$ws.Range("K41").Value = 15000
$ws.Range("M41").Value = -14562
# Row 93
$ws.Range("H93").Value = 2280.2727
$ws.Range("I93").Value = 1280.2858
$ws.Range("J93").Value = 4030.25
$ws.Range("K93").Value = 1280.2858
$ws.Range("L93").Value = 4030.25
$ws.Range("M93").Value = -32.28580000000011
$ws.Range("N93").Value = -6526.25
# Row 94
$ws.Range("H94").Value = 75000
$ws.Range("J94").Value = 75000
$ws.Range("L94").Value = 75000
$ws.Range("N94").Value = -76352
# Row 100
$ws.Range("H100").Value = 4708.727
$ws.Range("I100").Value = 3074.25
$ws.Range("K100").Value = 3074.25
$ws.Range("M100").Value = -2533.25
# Row 126
$ws.Range("H126").Value = 24256.75
$ws.Range("I126").Value = 19426.584
$ws.Range("K126").Value = 58279.75199999999
$ws.Range("M126").Value = -55809.75199999999
# Row 132
$ws.Range("H132").Value = 19620.977
$ws.Range("I132").Value = 23065.03
$ws.Range("K132").Value = 69195.09
$ws.Range("M132").Value = -66665.09

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 45
$ws.Range("H45").Value = 26182.5
$ws.Range("I45").Value = 58275
$ws.Range("J45").Value = 18159.375
$ws.Range("K45").Value = 58275
$ws.Range("L45").Value = 18159.375
$ws.Range("M45").Value = -57784
$ws.Range("N45").Value = -19141.375
# Row 124
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("M124").ClearContents()
$ws.Range("N124").ClearContents()
# Row 132
$ws.Range("H132").Value = 2502.2856
$ws.Range("I132").Value = 2540.6316
$ws.Range("K132").Value = 7621.8948
$ws.Range("M132").Value = -5091.8948
# Row 136
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
# Row 138
$ws.Range("H138").Value = 117500
$ws.Range("J138").Value = 145000
$ws.Range("L138").Value = 145000
$ws.Range("N138").Value = -155280
